$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.855.81'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '1.666.58'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '215.20'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").Value = '0.535'
$ws.Range("E6").Value = '  +5.17%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.252'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("E10").Value = '  +3.00%  '
$ws.Range("D12").Value = '1.902.38'
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = '1.656.43'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("D17").Value = '26.865.69'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '231.80'
$ws.Range("E18").Value = '  -3.83%  '
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  -1.98%  '
$ws.Range("E24").Value = '  -1.69%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").Value = '7.11'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("E27").Value = '  +1.43%  '
$ws.Range("D28").Value = '15.89'
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("D33").Value = '1.467.80'
$ws.Range("E33").Value = '  -4.06%  '
$ws.Range("D34").Value = '3.15'
$ws.Range("E34").Value = '  +3.18%  '
$ws.Range("E35").Value = '  +2.50%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '0.980'
$ws.Range("E43").Value = '  +6.71%  '
$ws.Range("D44").Value = '65.86'
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("D45").Value = '1.811.47'
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").Value = '0.778'
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").Value = '90.23'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.101'
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0508'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.54'
$ws.Range("E51").Value = '  +0.17%  '
